# This workbook's rows 67, 68, 70, 71, 72, 73 hold per-record species
# observation data (row 69 is untouched). The update cyclically rotates
# the record-specific fields among those six rows:
#   67 -> 68 -> 71 -> 72 -> 73 -> 70 -> 67
# i.e. the content that used to be in row 70 is now in row 67, the
# content that used to be in row 67 is now in row 68, and so on.
# Columns A,B,D,E,F,G,H,M,P,Q,R carry the record-specific data that
# moves; all other columns in this row range are identical across the
# six rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for each row after the rotation is applied.
$rowData = @{
    67 = @{ A = 111881310; B = 89425;  D = "NT"; E = 5442;   F = "Tallticka";             G = "Porodaedalea pini";               H = "(Brot.) Murrill";               M = $null;         P = "Valforsen, Ång" }
    68 = @{ A = 111870990; B = 90666;  D = "LC"; E = 4364;   F = "Dropptaggsvamp";         G = "Hydnellum ferrugineum";           H = "(Fr.:Fr.) P. Karst.";            M = $null;         P = "Valforsen (Valforsen), Ång" }
    70 = @{ A = 111870127; B = 89405;  D = "NT"; E = 1202;   F = "Ullticka";               G = "Phellinidium ferrugineofuscum";   H = "(P.Karst.) Fiasson & Niemelä";   M = $null;         P = "Valforsen (Valforsen), Ång" }
    71 = @{ A = 111881322; B = 56414;  D = "NT"; E = 100049; F = "Spillkråka";             G = "Dryocopus martius";               H = "(Linnaeus, 1758)";               M = "gammalt bo";  P = "Valforsen, Ång" }
    72 = @{ A = 111871585; B = 89405;  D = "NT"; E = 1202;   F = "Ullticka";               G = "Phellinidium ferrugineofuscum";   H = "(P.Karst.) Fiasson & Niemelä";   M = $null;         P = "Valforsen (Valforsen), Ång" }
    73 = @{ A = 111870830; B = 90678;  D = "LC"; E = 4366;   F = "Skarp dropptaggsvamp";   G = "Hydnellum peckii";                H = "Banker";                         M = $null;         P = "Valforsen (Valforsen), Ång" }
}

$qrData = @{
    67 = @{ Q = 590738.9206925276; R = 7040524.002523924 }
    68 = @{ Q = 590569.8478412227; R = 7040376.109235858 }
    70 = @{ Q = 590710.4131779457; R = 7040581.765558361 }
    71 = @{ Q = 590615.1562677342; R = 7040278.573758457 }
    72 = @{ Q = 590630.2636057099; R = 7040266.929520278 }
    73 = @{ Q = 590558.4251677697; R = 7040399.931061053 }
}

foreach ($r in @(67, 68, 70, 71, 72, 73)) {
    $d = $rowData[$r]
    $ws.Range("A$r").Value = $d.A
    $ws.Range("B$r").Value = $d.B
    $ws.Range("D$r").Value = $d.D
    $ws.Range("E$r").Value = $d.E
    $ws.Range("F$r").Value = $d.F
    $ws.Range("G$r").Value = $d.G
    $ws.Range("H$r").Value = $d.H
    if ($d.M) {
        $ws.Range("M$r").Value = $d.M
    } else {
        $ws.Range("M$r").ClearContents()
    }
    $ws.Range("P$r").Value = $d.P

    $qr = $qrData[$r]
    $ws.Range("Q$r").Value = $qr.Q
    $ws.Range("R$r").Value = $qr.R
}
